# BJ-999: Change signer in acts
#
# Update the signatory block (column A, rows 41-44) on the active sheet so
# the act is signed by the chief accountant instead of the chairman of the
# board:
#   A41  "Председатель Правления"            -> "Главный бухгалтер"
#   A42  "Губайдулин Т.Ф.,"                  -> "Кахно А.В.,"
#   A43  "действующий на основании Устава"   -> "действующая на основании"
#   A44  (blank)                             -> "Доверенности N 40 от 08.09.2020"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A41").Value = "Главный бухгалтер"
$ws.Range("A42").Value = "Кахно А.В.,"
$ws.Range("A43").Value = "действующая на основании"
$ws.Range("A44").Value = "Доверенности N 40 от 08.09.2020"

# Reflect the updated cursor/selection position recorded in the sheet view.
$ws.Range("C40").Select()
